$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New / updated body text for the e053d / e053e / e053f entries ---
$infantryResolutionText = "<Bold>e053d Main Gun Fire Against Infantry - Resolution</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nFor each hit scored against a target, consult the  `n<InlineUIContainer><Button Content='To Kill Infantry' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n Table to determine if the target is knocked out (KO'ed) using these modifier:`n<LineBreak/><LineBreak/>"
$vehicleResolutionText = "<Bold>e053e Main Gun Fire Against Vehicle - Resolution</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nFor each hit scored against a target, consult the correct To Kill Table to determine if the target is knocked out (KO'ed)."
$mgSelectTargetText = "<Bold>e053f Machine Gun (MG) Firing - Select Target</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nClick either the highlighted zone or an spotted target or <InlineUIContainer><Button Content='Skip MG' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> . `n<LineBreak/><LineBreak/>`n                                            <InlineUIContainer><Image Name='Continue53c' Height='100' Width='100'></Image></InlineUIContainer>"

# --- Insert a new row before row 70, pushing the evening-debrief rows (and everything below) down by one ---
$ws.Rows.Item(70).Insert()

# --- Row 68: "e053d" now documents the Infantry resolution (new, longer text) ---
$ws.Range("A68").Value = "e053d"
$ws.Range("B68").Value = $infantryResolutionText
$ws.Rows.Item(68).RowHeight = 115.2

# --- Row 69: "e053e" now documents the Vehicle resolution (former e053d body text, retitled) ---
$ws.Range("A69").Value = "e053e"
$ws.Range("B69").Value = $vehicleResolutionText
$ws.Rows.Item(69).RowHeight = 72

# --- Row 70 (newly inserted): "e053f" takes over the MG "Select Target" body text (former e053e body text, retitled) ---
$ws.Range("A70").Value = "e053f"
$ws.Range("B70").Value = $mgSelectTargetText
$ws.Rows.Item(70).RowHeight = 100.8

# --- Restore the saved selection/scroll state ---
$ws.Range("B66").Select()
